# Regenerate orders with updated distance/size codes.
#   Distances: D64 -> D69, D80 -> D86, D51 -> D55
#   Sizes:     S30 -> S31 (S20 / S25 unchanged)
# These substitutions are applied to every string-valued cell in the
# used range (condition labels, left/right filenames, and the
# Distance/Size lookup columns all encode the same tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            $nv = $v.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
